$d = $word.ActiveDocument

# Update the date heading (first paragraph)
$d.Paragraphs.Item(1).Range.Text = "2026-01-13 Tuesday"

# Update each answer cell in the 20x5 practice table
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "13+31=44"  # was "44+32=76"
$t.Cell(1, 2).Range.Text = "7+1=8"  # was "52-25=27"
$t.Cell(1, 3).Range.Text = "46-29=17"  # was "11+73=84"
$t.Cell(1, 4).Range.Text = "88-61=27"  # was "36-29=7"
$t.Cell(1, 5).Range.Text = "20+76=96"  # was "89-45=44"
$t.Cell(2, 1).Range.Text = "57-32=25"  # was "23-8=15"
$t.Cell(2, 2).Range.Text = "27+39=66"  # was "63+4=67"
$t.Cell(2, 3).Range.Text = "21+54=75"  # was "66+8=74"
$t.Cell(2, 4).Range.Text = "83-2=81"  # was "27-1=26"
$t.Cell(2, 5).Range.Text = "25+68=93"  # was "69-30=39"
$t.Cell(3, 1).Range.Text = "27-5=22"  # was "3+3=6"
$t.Cell(3, 2).Range.Text = "91-43=48"  # was "72-16=56"
$t.Cell(3, 3).Range.Text = "67-31=36"  # was "67+11=78"
$t.Cell(3, 4).Range.Text = "87-5=82"  # was "80-75=5"
$t.Cell(3, 5).Range.Text = "79-30=49"  # was "7+13=20"
$t.Cell(4, 1).Range.Text = "27+65=92"  # was "78+15=93"
$t.Cell(4, 2).Range.Text = "43-31=12"  # was "88-13=75"
$t.Cell(4, 3).Range.Text = "23-4=19"  # was "99-4=95"
$t.Cell(4, 4).Range.Text = "81-5=76"  # was "29-3=26"
$t.Cell(4, 5).Range.Text = "36-36=0"  # was "40+19=59"
$t.Cell(5, 1).Range.Text = "69+21=90"  # was "1+80=81"
$t.Cell(5, 2).Range.Text = "53+17=70"  # was "50+9=59"
$t.Cell(5, 3).Range.Text = "46-7=39"  # was "36+37=73"
$t.Cell(5, 4).Range.Text = "51+39=90"  # was "18-11=7"
$t.Cell(5, 5).Range.Text = "34+26=60"  # was "87-58=29"
$t.Cell(6, 1).Range.Text = "37-4=33"  # was "67-2=65"
$t.Cell(6, 2).Range.Text = "25-9=16"  # was "72+20=92"
$t.Cell(6, 3).Range.Text = "44+38=82"  # was "34+25=59"
$t.Cell(6, 4).Range.Text = "46+12=58"  # was "3+50=53"
$t.Cell(6, 5).Range.Text = "89+1=90"  # was "41-39=2"
$t.Cell(7, 1).Range.Text = "23+8=31"  # was "75+0=75"
$t.Cell(7, 2).Range.Text = "44+50=94"  # was "53+19=72"
$t.Cell(7, 3).Range.Text = "18+51=69"  # was "45-6=39"
$t.Cell(7, 4).Range.Text = "9+9=18"  # was "45-12=33"
$t.Cell(7, 5).Range.Text = "92-55=37"  # was "47+4=51"
$t.Cell(8, 1).Range.Text = "39-1=38"  # was "56+13=69"
$t.Cell(8, 2).Range.Text = "24+48=72"  # was "34-4=30"
$t.Cell(8, 3).Range.Text = "24+72=96"  # was "9+12=21"
$t.Cell(8, 4).Range.Text = "26+16=42"  # was "98-54=44"
$t.Cell(8, 5).Range.Text = "67-36=31"  # was "97-59=38"
$t.Cell(9, 1).Range.Text = "83-5=78"  # was "76-49=27"
$t.Cell(9, 2).Range.Text = "87-69=18"  # was "85-16=69"
$t.Cell(9, 3).Range.Text = "0+71=71"  # was "90-14=76"
$t.Cell(9, 4).Range.Text = "84-83=1"  # was "43-38=5"
$t.Cell(9, 5).Range.Text = "56-14=42"  # was "46-46=0"
$t.Cell(10, 1).Range.Text = "68-6=62"  # was "21-9=12"
$t.Cell(10, 2).Range.Text = "50-35=15"  # was "28+48=76"
$t.Cell(10, 3).Range.Text = "12+73=85"  # was "60-41=19"
$t.Cell(10, 4).Range.Text = "90-21=69"  # was "48+13=61"
$t.Cell(10, 5).Range.Text = "84-25=59"  # was "28+51=79"
$t.Cell(11, 1).Range.Text = "21-14=7"  # was "92-28=64"
$t.Cell(11, 2).Range.Text = "48-44=4"  # was "4+86=90"
$t.Cell(11, 3).Range.Text = "64-20=44"  # was "68+16=84"
$t.Cell(11, 4).Range.Text = "93-64=29"  # was "99-5=94"
$t.Cell(11, 5).Range.Text = "72-62=10"  # was "18+48=66"
$t.Cell(12, 1).Range.Text = "81-65=16"  # was "14+30=44"
$t.Cell(12, 2).Range.Text = "3+44=47"  # was "1+47=48"
$t.Cell(12, 3).Range.Text = "75-49=26"  # was "16+5=21"
$t.Cell(12, 4).Range.Text = "36+26=62"  # was "34+5=39"
$t.Cell(12, 5).Range.Text = "1+50=51"  # was "86-11=75"
$t.Cell(13, 1).Range.Text = "46+18=64"  # was "43+55=98"
$t.Cell(13, 2).Range.Text = "97-5=92"  # was "90-22=68"
$t.Cell(13, 3).Range.Text = "43+50=93"  # was "71+10=81"
$t.Cell(13, 4).Range.Text = "89+0=89"  # was "38-32=6"
$t.Cell(13, 5).Range.Text = "71-69=2"  # was "6+92=98"
$t.Cell(14, 1).Range.Text = "96-15=81"  # was "84-44=40"
$t.Cell(14, 2).Range.Text = "29+54=83"  # was "23-3=20"
$t.Cell(14, 3).Range.Text = "17-15=2"  # was "34-1=33"
$t.Cell(14, 4).Range.Text = "9+27=36"  # was "28+6=34"
$t.Cell(14, 5).Range.Text = "93-77=16"  # was "10+60=70"
$t.Cell(15, 1).Range.Text = "39+51=90"  # was "91-53=38"
$t.Cell(15, 2).Range.Text = "88+8=96"  # was "39+3=42"
$t.Cell(15, 3).Range.Text = "14+3=17"  # was "59-29=30"
$t.Cell(15, 4).Range.Text = "16+35=51"  # was "74-7=67"
$t.Cell(15, 5).Range.Text = "79-71=8"  # was "40+14=54"
$t.Cell(16, 1).Range.Text = "79-20=59"  # was "23+65=88"
$t.Cell(16, 2).Range.Text = "74-67=7"  # was "66-62=4"
$t.Cell(16, 3).Range.Text = "57+8=65"  # was "1+44=45"
$t.Cell(16, 4).Range.Text = "71-22=49"  # was "8+25=33"
$t.Cell(16, 5).Range.Text = "84-29=55"  # was "97-51=46"
$t.Cell(17, 1).Range.Text = "18+75=93"  # was "19+11=30"
$t.Cell(17, 2).Range.Text = "37-33=4"  # was "21+71=92"
$t.Cell(17, 3).Range.Text = "18+42=60"  # was "70+27=97"
$t.Cell(17, 4).Range.Text = "85+2=87"  # was "44-28=16"
$t.Cell(17, 5).Range.Text = "31-16=15"  # was "59-29=30"
$t.Cell(18, 1).Range.Text = "85-34=51"  # was "12+77=89"
$t.Cell(18, 2).Range.Text = "75-26=49"  # was "74-8=66"
$t.Cell(18, 3).Range.Text = "88-79=9"  # was "41-1=40"
$t.Cell(18, 4).Range.Text = "90-11=79"  # was "47+36=83"
$t.Cell(18, 5).Range.Text = "59-17=42"  # was "88-28=60"
$t.Cell(19, 1).Range.Text = "67-58=9"  # was "38+17=55"
$t.Cell(19, 2).Range.Text = "9+30=39"  # was "79-79=0"
$t.Cell(19, 3).Range.Text = "99-46=53"  # was "66+15=81"
$t.Cell(19, 4).Range.Text = "47+9=56"  # was "44+49=93"
$t.Cell(19, 5).Range.Text = "69-19=50"  # was "74+24=98"
$t.Cell(20, 1).Range.Text = "87-41=46"  # was "50-16=34"
$t.Cell(20, 2).Range.Text = "57+22=79"  # was "81-68=13"
$t.Cell(20, 3).Range.Text = "18+28=46"  # was "51+24=75"
$t.Cell(20, 4).Range.Text = "66+21=87"  # was "12+3=15"
$t.Cell(20, 5).Range.Text = "91-79=12"  # was "21-18=3"
